$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells hold plain text that often *looks* numeric
# ("5.03", "2.00", "0.0000220", ...). Excel's COM Value setter auto-detects
# numbers from such strings, which would silently rewrite them (dropping
# trailing zeros / exact formatting). Force Text number format while writing,
# then clear the format again so the cell style index is unchanged (matches
# the original workbook, which has no explicit style on these cells).
$dCells = $ws.Range("D2","D3","D5","D6","D8","D9","D10","D11","D13","D14","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D29","D31","D32","D33","D34","D35","D38","D39","D42","D43","D44","D46","D50","D51")
$dCells.NumberFormat = "@"

$ws.Range("D2").Value = "60.087.28"
$ws.Range("D3").Value = "2.993.82"
$ws.Range("D5").Value = "569.89"
$ws.Range("D6").Value = "124.92"
$ws.Range("D8").Value = "2.987.09"
$ws.Range("D9").Value = "0.503"
$ws.Range("D10").Value = "0.131"
$ws.Range("D11").Value = "5.03"
$ws.Range("D13").Value = "0.0000220"
$ws.Range("D14").Value = "32.39"
$ws.Range("D16").Value = "3.489.78"
$ws.Range("D17").Value = "2.996.80"
$ws.Range("D18").Value = "60.089.20"
$ws.Range("D19").Value = "6.53"
$ws.Range("D20").Value = "427.34"
$ws.Range("D21").Value = "13.13"
$ws.Range("D22").Value = "0.670"
$ws.Range("D23").Value = "7.06"
$ws.Range("D24").Value = "12.87"
$ws.Range("D25").Value = "79.45"
$ws.Range("D29").Value = "1.95"
$ws.Range("D31").Value = "25.22"
$ws.Range("D32").Value = "6.07"
$ws.Range("D33").Value = "0.0950"
$ws.Range("D34").Value = "5.58"
$ws.Range("D35").Value = "0.928"
$ws.Range("D38").Value = "8.48"
$ws.Range("D39").Value = "0.0₃0656"
$ws.Range("D42").Value = "370.21"
$ws.Range("D43").Value = "2.666.19"
$ws.Range("D44").Value = "2.43"
$ws.Range("D46").Value = "120.93"
$ws.Range("D50").Value = "23.23"
$ws.Range("D51").Value = "2.00"

$dCells.ClearFormats()

# Other columns (Coin name, Link URL, Volume%) never look like plain numbers
# (URLs/names are alphanumeric; percentages keep their surrounding spaces and
# "%" suffix), so a direct .Value assignment keeps them as text safely.
$ws.Range("E2").Value = "  -5.07%  "
$ws.Range("E3").Value = "  -5.48%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  -4.92%  "
$ws.Range("E6").Value = "  -7.77%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("E8").Value = "  -5.69%  "
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("E10").Value = "  -7.82%  "
$ws.Range("E11").Value = "  -4.97%  "
$ws.Range("E12").Value = "  -2.72%  "
$ws.Range("E13").Value = "  -8.11%  "
$ws.Range("E14").Value = "  -6.82%  "
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("E16").Value = "  -5.42%  "
$ws.Range("E17").Value = "  -5.34%  "
$ws.Range("E18").Value = "  -5.02%  "
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("E20").Value = "  -7.30%  "
$ws.Range("E21").Value = "  -6.22%  "
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("E23").Value = "  -7.71%  "
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("E25").Value = "  -4.34%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  -6.53%  "
$ws.Range("E29").Value = "  -5.63%  "
$ws.Range("E30").Value = "  -6.98%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E31").Value = "  -7.13%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E32").Value = "  -10.62%  "
$ws.Range("E33").Value = "  -5.36%  "
$ws.Range("E34").Value = "  -4.91%  "
$ws.Range("E35").Value = "  -8.90%  "
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("E37").Value = "  -16.55%  "
$ws.Range("E38").Value = "  +4.14%  "
$ws.Range("E39").Value = "  -10.58%  "
$ws.Range("E40").Value = "  -9.08%  "
$ws.Range("E41").Value = "  -4.31%  "
$ws.Range("E42").Value = "  -5.69%  "
$ws.Range("E43").Value = "  -4.45%  "
$ws.Range("E44").Value = "  -7.33%  "
$ws.Range("E46").Value = "  -4.76%  "
$ws.Range("E47").Value = "  -7.06%  "
$ws.Range("E48").Value = "  -6.35%  "
$ws.Range("E49").Value = "  -3.53%  "
$ws.Range("E50").Value = "  -7.29%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("E51").Value = "  -7.12%  "
